$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Student model update for the new semester: add a new (currently blank/placeholder)
# row to the gridnodes sheet by writing a single space into H6.
$ws.Range("H6").Value = " "

# Reflect the new active selection from the saved workbook view.
$ws.Range("J12").Select()
